# Add product cards export rows (base_id/quote records) for 2025-11-19 and
# 2025-11-20 to the records sheet, rows 18-22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 18; A = "20251119-002"; B = "2025-11-19"; C = "q"; D = "QUO-20251119-001"; E = 115;   H = "Abu Dhabi - Al Shamkha" },
    @{ Row = 19; A = "20251120-003"; B = "2025-11-20"; C = "q"; D = "QUO-20251120-005"; E = 985;   H = "Abu Dhabi - Al Shamkha" },
    @{ Row = 20; A = "20251120-004"; B = "2025-11-20"; C = "q"; D = "QUO-20251120-003"; E = 425;   H = "Abu Dhabi - Al Shamkha" },
    @{ Row = 21; A = "20251120-004"; B = "2025-11-20"; C = "q"; D = "QUO-20251120-001"; E = 14030; H = "Abu Dhabi - Al Shamkha" },
    @{ Row = 22; A = "20251120-005"; B = "2025-11-20"; C = "q"; D = "QUO-20251120-033"; E = 14030; H = "Abu Dhabi - Al Shamkha" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Columns A, C, D, F, G, H, I are plain text. Column B looks like a date
    # ("YYYY-MM-DD") and column A/D contain digit-leading tokens, so force a
    # text number format before assignment to stop Excel's automatic
    # type-sniffing (dates / numbers), then clear the format again so no
    # style index is left behind on the cell (matches the source rows,
    # which carry no explicit style).
    $ws.Range("A$rowNum").NumberFormat = "@"
    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("A$rowNum").ClearFormats()

    $ws.Range("B$rowNum").NumberFormat = "@"
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("B$rowNum").ClearFormats()

    $ws.Range("C$rowNum").NumberFormat = "@"
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("C$rowNum").ClearFormats()

    $ws.Range("D$rowNum").NumberFormat = "@"
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("D$rowNum").ClearFormats()

    # amount column is numeric
    $ws.Range("E$rowNum").Value = $r.E

    # client_name / phone / note are empty text cells for these rows
    $ws.Range("F$rowNum").NumberFormat = "@"
    $ws.Range("F$rowNum").Value = ""
    $ws.Range("F$rowNum").ClearFormats()

    $ws.Range("G$rowNum").NumberFormat = "@"
    $ws.Range("G$rowNum").Value = ""
    $ws.Range("G$rowNum").ClearFormats()

    $ws.Range("H$rowNum").NumberFormat = "@"
    $ws.Range("H$rowNum").Value = $r.H
    $ws.Range("H$rowNum").ClearFormats()

    $ws.Range("I$rowNum").NumberFormat = "@"
    $ws.Range("I$rowNum").Value = ""
    $ws.Range("I$rowNum").ClearFormats()
}
